try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # A new weekly price record was reported for "Provincia de Limarí" that
    # belongs chronologically right before the existing row 125 entry, so a
    # row is inserted at 125 and the following rows (old 125-131) shift down
    # to 126-132.
    $ws.Rows(125).Insert()

    # Populate the newly inserted row 125 with the new record's data.
    $ws.Cells.Item(125, 1).Value = 8
    $ws.Cells.Item(125, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item(125, 3).Value = "Coquimbo"
    $ws.Cells.Item(125, 4).Value = 45021
    $ws.Cells.Item(125, 5).Value = 4
    $ws.Cells.Item(125, 6).Value = 100112030
    $ws.Cells.Item(125, 7).Value = "Poroto granado"
    $ws.Cells.Item(125, 8).Value = "Sin especificar"
    $ws.Cells.Item(125, 9).Value = "Primera"
    $ws.Cells.Item(125, 10).Value = 400
    $ws.Cells.Item(125, 11).Value = 29000
    $ws.Cells.Item(125, 12).Value = 30000
    $ws.Cells.Item(125, 13).Value = 29500
    $ws.Cells.Item(125, 14).Value = "`$/malla 25 kilos"
    $ws.Cells.Item(125, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item(125, 16).Value = 1180
    $ws.Cells.Item(125, 17).Value = 25
    $ws.Cells.Item(125, 18).Value = "Hortaliza"
} catch {
    Write-Output "ERROR: $_"
}
